$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.112.23"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "3.601.80"
$ws.Range("E3").Value = "  +2.45%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.40"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "196.47"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -1.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.648"
$ws.Range("E10").Value = "  -0.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.74"
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("E12").Value = "  +0.90%  "
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").Value = "4.177.00"
$ws.Range("E14").Value = "  +2.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.06"
$ws.Range("E15").Value = "  +3.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "594.20"
$ws.Range("E16").Value = "  -1.54%  "
$ws.Range("D17").Value = "70.287.94"
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.10"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").Value = "3.603.92"
$ws.Range("E19").Value = "  +2.53%  "
$ws.Range("E20").Value = "  +1.36%  "
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.73"
$ws.Range("E22").Value = "  -2.64%  "
$ws.Range("E23").Value = "  -1.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "101.52"
$ws.Range("E24").Value = "  -2.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.60"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.01"
$ws.Range("E26").Value = "  -1.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.73"
$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.57"
$ws.Range("E28").Value = "  -0.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.74"
$ws.Range("E29").Value = "  +0.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.67"
$ws.Range("E30").Value = "  +3.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.14"
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.28"
$ws.Range("E32").Value = "  -3.49%  "
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.24"
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("D35").Value = "0.0₃0879"
$ws.Range("E35").Value = "  +7.51%  "
$ws.Range("D36").Value = "3.921.89"
$ws.Range("E36").Value = "  +4.80%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "525.54"
$ws.Range("E37").Value = "  +7.04%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.10"
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.81"
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("E42").Value = "  -1.84%  "
$ws.Range("E43").Value = "  -1.98%  "
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.45"
$ws.Range("E45").Value = "  +3.21%  "
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.61"
$ws.Range("E48").Value = "  -0.52%  "
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000249"
$ws.Range("E50").Value = "  +2.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.31"
$ws.Range("E51").Value = "  +1.70%  "
